$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.913
$ws.Range("A6").Value = -21.14
$ws.Range("A7").Value = -21.089
$ws.Range("A8").Value = -21.018
$ws.Range("A16").Value = -20.727
$ws.Range("A20").Value = -22.138
$ws.Range("A21").Value = -21.14
